# JOHNPAUL_20250805_cycle1.xlsx - sync repo with latest training pipeline, results reorg, and artifacts
#
# Net data edit: the vehicle Count for the JP_W / Bus row (row 23) in
# Raw_Annotations was cleared out (previously 4, now blank). Every other
# changed cell in the workbook (Raw_Annotations!G23/I23/J23/K23 and the
# Aggregates!row-5 SUMIFS block) is a formula that depends on that cell, so
# Excel recalculates them automatically once the value is cleared.
#
# The workbook also no longer carries the reference photo that used to be
# anchored on the Raw_Annotations sheet (xl/drawings/drawing1.xml), so we
# remove that picture too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Raw_Annotations")

# --- Remove the embedded reference picture (image1.png) from Raw_Annotations ---
for ($i = $ws.Shapes.Count; $i -ge 1; $i--) {
    $ws.Shapes.Item($i).Delete()
}

# --- Clear the vehicle count for JP_W / Bus (row 23) ---
$ws.Range("D23").ClearContents()

# Recalculate so every dependent formula (G23/I23/J23/K23 on this sheet and
# the JP_W aggregate row on the Aggregates sheet) gets a fresh cached value.
$excel.CalculateFull()
